# Auto-generated edit script: update Sargatanas_Profits market-value columns
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: H2=47.8, I2=32.307693, K2=32.307693, M2=80.692307
$ws.Range("H2").Value = 47.8
$ws.Range("I2").Value = 32.307693
$ws.Range("K2").Value = 32.307693
$ws.Range("M2").Value = 80.692307

# Row 17: H17=1281.8077, J17=1301.08, L17=3903.24, N17=-4239.24
$ws.Range("H17").Value = 1281.8077
$ws.Range("J17").Value = 1301.08
$ws.Range("L17").Value = 3903.24
$ws.Range("N17").Value = -4239.24

# Row 40: H40=5212.5713, I40=2997.6667, K40=2997.6667, M40=-2822.6667
$ws.Range("H40").Value = 5212.5713
$ws.Range("I40").Value = 2997.6667
$ws.Range("K40").Value = 2997.6667
$ws.Range("M40").Value = -2822.6667

# Row 98: H98=6416.091, I98=6416.091, K98=6416.091, M98=-4918.091
$ws.Range("H98").Value = 6416.091
$ws.Range("I98").Value = 6416.091
$ws.Range("K98").Value = 6416.091
$ws.Range("M98").Value = -4918.091

# Row 122: H122=6416.091, I122=6416.091, K122=19248.273, M122=-16798.273
$ws.Range("H122").Value = 6416.091
$ws.Range("I122").Value = 6416.091
$ws.Range("K122").Value = 19248.273
$ws.Range("M122").Value = -16798.273

# Row 135: H135=455115.8, I135=769974, J135=320.66666, K135=6929766, L135=2885.99994, M135=-6927231, N135=-7955.99994
$ws.Range("H135").Value = 455115.8
$ws.Range("I135").Value = 769974
$ws.Range("J135").Value = 320.66666
$ws.Range("K135").Value = 6929766
$ws.Range("L135").Value = 2885.99994
$ws.Range("M135").Value = -6927231
$ws.Range("N135").Value = -7955.99994

# Row 137: H137=2604.2258, I137=2195.7646, K137=6587.293799999999, M137=-4037.293799999999
$ws.Range("H137").Value = 2604.2258
$ws.Range("I137").Value = 2195.7646
$ws.Range("K137").Value = 6587.293799999999
$ws.Range("M137").Value = -4037.293799999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=1566338.4, I32=1692623.8, J32=8818, K32=1692623.8, L32=8818, M32=-1692336.8, N32=-9392
$ws.Range("H32").Value = 1566338.4
$ws.Range("I32").Value = 1692623.8
$ws.Range("J32").Value = 8818
$ws.Range("K32").Value = 1692623.8
$ws.Range("L32").Value = 8818
$ws.Range("M32").Value = -1692336.8
$ws.Range("N32").Value = -9392

# Row 61: H61=5585.608, I61=3177.9707, K61=3177.9707, M61=-2965.9707
$ws.Range("H61").Value = 5585.608
$ws.Range("I61").Value = 3177.9707
$ws.Range("K61").Value = 3177.9707
$ws.Range("M61").Value = -2965.9707

# Row 102: H102=2033, I102=2099.5, K102=2099.5, M102=-477.5
$ws.Range("H102").Value = 2033
$ws.Range("I102").Value = 2099.5
$ws.Range("K102").Value = 2099.5
$ws.Range("M102").Value = -477.5

# Row 136: H136=5585.608, I136=3177.9707, K136=9533.9121, M136=-6983.9121
$ws.Range("H136").Value = 5585.608
$ws.Range("I136").Value = 3177.9707
$ws.Range("K136").Value = 9533.9121
$ws.Range("M136").Value = -6983.9121

$ws = $wb.Worksheets.Item("BSM")
# Row 86: H86=35717300, I86=1752.25, J86=83338030, K86=1752.25, L86=83338030, M86=-629.25, N86=-83340276
$ws.Range("H86").Value = 35717300
$ws.Range("I86").Value = 1752.25
$ws.Range("J86").Value = 83338030
$ws.Range("K86").Value = 1752.25
$ws.Range("L86").Value = 83338030
$ws.Range("M86").Value = -629.25
$ws.Range("N86").Value = -83340276

# Row 89: H89=35717300, I89=1752.25, J89=83338030, K89=8761.25, L89=416690150, M89=-3145.25, N89=-416701382
$ws.Range("H89").Value = 35717300
$ws.Range("I89").Value = 1752.25
$ws.Range("J89").Value = 83338030
$ws.Range("K89").Value = 8761.25
$ws.Range("L89").Value = 416690150
$ws.Range("M89").Value = -3145.25
$ws.Range("N89").Value = -416701382

# Row 99: H99=2757624.2, I99=2545.5, K99=2545.5, M99=-1047.5
$ws.Range("H99").Value = 2757624.2
$ws.Range("I99").Value = 2545.5
$ws.Range("K99").Value = 2545.5
$ws.Range("M99").Value = -1047.5

# Row 105: H105=2914.8462, I105=2525.9473, J105=3970.4285, K105=2525.9473, L105=3970.4285, M105=-778.9472999999998, N105=-7464.4285
$ws.Range("H105").Value = 2914.8462
$ws.Range("I105").Value = 2525.9473
$ws.Range("J105").Value = 3970.4285
$ws.Range("K105").Value = 2525.9473
$ws.Range("L105").Value = 3970.4285
$ws.Range("M105").Value = -778.9472999999998
$ws.Range("N105").Value = -7464.4285

# Row 107: H107=86544504, I107=112504360, J107=11663, K107=112504360, L107=11663, M107=-112502440, N107=-15503
$ws.Range("H107").Value = 86544504
$ws.Range("I107").Value = 112504360
$ws.Range("J107").Value = 11663
$ws.Range("K107").Value = 112504360
$ws.Range("L107").Value = 11663
$ws.Range("M107").Value = -112502440
$ws.Range("N107").Value = -15503

# Row 128: H128=2348, I128=2348, K128=7044, M128=-4554
$ws.Range("H128").Value = 2348
$ws.Range("I128").Value = 2348
$ws.Range("K128").Value = 7044
$ws.Range("M128").Value = -4554

# Row 134: H134=5122.585, I134=1887.8572, K134=5663.571599999999, M134=-3128.571599999999
$ws.Range("H134").Value = 5122.585
$ws.Range("I134").Value = 1887.8572
$ws.Range("K134").Value = 5663.571599999999
$ws.Range("M134").Value = -3128.571599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 108: H108=34998.5, J108=59376, L108=59376, N108=-67056
$ws.Range("H108").Value = 34998.5
$ws.Range("J108").Value = 59376
$ws.Range("L108").Value = 59376
$ws.Range("N108").Value = -67056

$ws = $wb.Worksheets.Item("CUL")
# Row 107: H107=20000354, J107=28571762, L107=85715286, N107=-85719126
$ws.Range("H107").Value = 20000354
$ws.Range("J107").Value = 28571762
$ws.Range("L107").Value = 85715286
$ws.Range("N107").Value = -85719126

# Row 127: H127=1000, J127=1000, L127=3000, N127=-12920
$ws.Range("H127").Value = 1000
$ws.Range("J127").Value = 1000
$ws.Range("L127").Value = 3000
$ws.Range("N127").Value = -12920

# Row 131: H131=1436.1, I131=670.375, K131=2011.125, M131=3028.875
$ws.Range("H131").Value = 1436.1
$ws.Range("I131").Value = 670.375
$ws.Range("K131").Value = 2011.125
$ws.Range("M131").Value = 3028.875

# Row 133: H133=4336, I133=998.5, K133=2995.5, M133=2064.5
$ws.Range("H133").Value = 4336
$ws.Range("I133").Value = 998.5
$ws.Range("K133").Value = 2995.5
$ws.Range("M133").Value = 2064.5

# Row 137: H137=45158.348, I137=1156.3077, K137=3468.9231, M137=1631.0769
$ws.Range("H137").Value = 45158.348
$ws.Range("I137").Value = 1156.3077
$ws.Range("K137").Value = 3468.9231
$ws.Range("M137").Value = 1631.0769

# Row 139: H139=6730.619, I139=2988.2222, J139=9537.416999999999, K139=8964.6666, L139=28612.251, M139=-3824.6666, N139=-38892.251
$ws.Range("H139").Value = 6730.619
$ws.Range("I139").Value = 2988.2222
$ws.Range("J139").Value = 9537.416999999999
$ws.Range("K139").Value = 8964.6666
$ws.Range("L139").Value = 28612.251
$ws.Range("M139").Value = -3824.6666
$ws.Range("N139").Value = -38892.251

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70=10131.85, I70=8955.666999999999, K70=8955.666999999999, M70=-8685.666999999999
$ws.Range("H70").Value = 10131.85
$ws.Range("I70").Value = 8955.666999999999
$ws.Range("K70").Value = 8955.666999999999
$ws.Range("M70").Value = -8685.666999999999

# Row 73: H73=10131.85, I73=8955.666999999999, K73=8955.666999999999, M73=-8019.666999999999
$ws.Range("H73").Value = 10131.85
$ws.Range("I73").Value = 8955.666999999999
$ws.Range("K73").Value = 8955.666999999999
$ws.Range("M73").Value = -8019.666999999999

# Row 113: H113=6548.353, I113=2908.5386, J113=8801.571, K113=2908.5386, L113=8801.571, M113=-738.5385999999999, N113=-13141.571
$ws.Range("H113").Value = 6548.353
$ws.Range("I113").Value = 2908.5386
$ws.Range("J113").Value = 8801.571
$ws.Range("K113").Value = 2908.5386
$ws.Range("L113").Value = 8801.571
$ws.Range("M113").Value = -738.5385999999999
$ws.Range("N113").Value = -13141.571

$ws = $wb.Worksheets.Item("LTW")
# Row 46: H46=2214.45, I46=472.77777, J46=3639.4546, K46=472.77777, L46=3639.4546, M46=-284.77777, N46=-4015.4546
$ws.Range("H46").Value = 2214.45
$ws.Range("I46").Value = 472.77777
$ws.Range("J46").Value = 3639.4546
$ws.Range("K46").Value = 472.77777
$ws.Range("L46").Value = 3639.4546
$ws.Range("M46").Value = -284.77777
$ws.Range("N46").Value = -4015.4546

# Row 93: H93=5005.478, I93=5794.2, J93=4398.769, K93=5794.2, L93=4398.769, M93=-4546.2, N93=-6894.769
$ws.Range("H93").Value = 5005.478
$ws.Range("I93").Value = 5794.2
$ws.Range("J93").Value = 4398.769
$ws.Range("K93").Value = 5794.2
$ws.Range("L93").Value = 4398.769
$ws.Range("M93").Value = -4546.2
$ws.Range("N93").Value = -6894.769

# Row 100: H100=3878.0625, I100=2754.4, K100=2754.4, M100=-2213.4
$ws.Range("H100").Value = 3878.0625
$ws.Range("I100").Value = 2754.4
$ws.Range("K100").Value = 2754.4
$ws.Range("M100").Value = -2213.4

$ws = $wb.Worksheets.Item("WVR")
# Row 4: H4=9756.9, I4=9285.714, J4=10856.333, K4=9285.714, L4=10856.333, M4=-9172.714, N4=-11082.333
$ws.Range("H4").Value = 9756.9
$ws.Range("I4").Value = 9285.714
$ws.Range("J4").Value = 10856.333
$ws.Range("K4").Value = 9285.714
$ws.Range("L4").Value = 10856.333
$ws.Range("M4").Value = -9172.714
$ws.Range("N4").Value = -11082.333

# Row 15: H15=24995.6, I15=24995.6, K15=24995.6, M15=-24707.6
$ws.Range("H15").Value = 24995.6
$ws.Range("I15").Value = 24995.6
$ws.Range("K15").Value = 24995.6
$ws.Range("M15").Value = -24707.6

# Row 100: H100=720.1739, J100=764.46155, L100=1528.9231, N100=-2610.9231
$ws.Range("H100").Value = 720.1739
$ws.Range("J100").Value = 764.46155
$ws.Range("L100").Value = 1528.9231
$ws.Range("N100").Value = -2610.9231

# Row 107: H107=27778580, I107=415.375, K107=1246.125, M107=673.875
$ws.Range("H107").Value = 27778580
$ws.Range("I107").Value = 415.375
$ws.Range("K107").Value = 1246.125
$ws.Range("M107").Value = 673.875

# Row 126: H126=2958.5925, I126=1898.25, J126=3806.8667, K126=5694.75, L126=11420.6001, M126=-3224.75, N126=-16360.6001
$ws.Range("H126").Value = 2958.5925
$ws.Range("I126").Value = 1898.25
$ws.Range("J126").Value = 3806.8667
$ws.Range("K126").Value = 5694.75
$ws.Range("L126").Value = 11420.6001
$ws.Range("M126").Value = -3224.75
$ws.Range("N126").Value = -16360.6001
